$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = (Get-Date -Year 2025 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("B2").Value = 73.28
$ws.Range("C2").Value = 62.28
$ws.Range("D2").Value = 41.87
$ws.Range("E2").Value = 39.39
$ws.Range("F2").Value = 34.63
$ws.Range("G2").Value = 42.09
$ws.Range("H2").Value = 47.43
$ws.Range("I2").Value = 64.20999999999999
$ws.Range("J2").Value = 58.06
$ws.Range("K2").Value = 16.44
$ws.Range("L2").Value = 1.64
$ws.Range("M2").Value = 0.03
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.28
$ws.Range("R2").Value = 14.1
$ws.Range("S2").Value = 58.81
$ws.Range("T2").Value = 68.76000000000001
$ws.Range("U2").Value = 73.51000000000001
$ws.Range("V2").Value = 79.15000000000001
$ws.Range("W2").Value = 72.65000000000001
$ws.Range("X2").Value = 65.19
$ws.Range("Y2").Value = 43.26
$ws.Range("Z2").Value = 39.88

$ws.Range("AB2").Value = 65.06
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 75.90000000000001
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 71.14
$ws.Range("AG2").Value = "3h-16h"
